$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new site/login/password entry
$ws.Range("A4").Value = "www.f"
$ws.Range("B4").Value = "dima"
$ws.Range("C4").Value = "LRvWxWSA"

# Row 5: new site/login/password entry (login reuses existing "DDfire" value)
$ws.Range("A5").Value = "www.dima"
$ws.Range("B5").Value = "DDfire"
$ws.Range("C5").Value = "zR{X7/9od7nF"
